$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "model" sheet: add an "raw" element (object / mime file) definition.
# A new header-row set of columns (C:H) is introduced, plus a third data
# row describing the new "raw" field (mimeUri/mimeType object).
# ---------------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")

$model.Range("C1").Value = "elementType"
$model.Range("D1").Value = "properties.uriFragment.type"
$model.Range("E1").Value = "properties.uriFragment.elementType"
$model.Range("F1").Value = "properties.contentType.type"
$model.Range("G1").Value = "properties.contentType.elementType"
$model.Range("H1").Value = "properties.contentType.default"

$model.Range("A3").Value = "object"
$model.Range("B3").Value = "raw"
$model.Range("C3").Value = "mimeUri"
$model.Range("D3").Value = "string"
$model.Range("E3").Value = "rowpath"
$model.Range("F3").Value = "string"
$model.Range("G3").Value = "mimeType"
$model.Range("H3").Value = "application/json"

# New column widths introduced alongside the new columns.
$model.Columns.Item(4).ColumnWidth = 24.998697916666668
$model.Columns.Item(5).ColumnWidth = 31.830729166666668
$model.Columns.Item(6).ColumnWidth = 33.166666666666664
$model.Columns.Item(7).ColumnWidth = 33.166666666666664
$model.Columns.Item(8).ColumnWidth = 28.330729166666668

$model.PageSetup.Orientation = 1

$model.Range("A3").Select()

# ---------------------------------------------------------------------------
# "properties" sheet: the "colOrder" array (Table/default/colOrder) now
# also lists the new raw/raw_contentType/raw_uriFragment columns.
# ---------------------------------------------------------------------------
$properties = $wb.Worksheets.Item("properties")
$properties.Range("E2").Value = '["address","address_image0_contentType","address_image0_uriFragment","comments","comments_image0_contentType","comments_image0_uriFragment","fri_chores","fri_chores_image0_contentType","fri_chores_image0_uriFragment","mon_chores","mon_chores_image0_contentType","mon_chores_image0_uriFragment","name","name_image0_contentType","name_image0_uriFragment","qrcode","qrcode_image0_contentType","qrcode_image0_uriFragment","raw","raw_contentType","raw_uriFragment","roomNum","roomNum_image0_contentType","roomNum_image0_uriFragment","sat_chores","sat_chores_image0_contentType","sat_chores_image0_uriFragment","scan_output_directory","stay","stay_image0_contentType","stay_image0_uriFragment","sun_chores","sun_chores_image0_contentType","sun_chores_image0_uriFragment","thurs_chores","thurs_chores_image0_contentType","thurs_chores_image0_uriFragment","tues_chores","tues_chores_image0_contentType","tues_chores_image0_uriFragment","wed_chores","wed_chores_image0_contentType","wed_chores_image0_uriFragment"]'

# ---------------------------------------------------------------------------
# Minor cosmetic re-selections / column-width touch-ups left behind by the
# resave (best effort; values come from the target column character widths).
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Columns.Item(2).ColumnWidth = 18.666666666666668
$survey.Columns.Item(4).ColumnWidth = 30.330729166666668
$survey.Range("B52").Select()

$choices = $wb.Worksheets.Item("choices")
$choices.Columns.Item(1).ColumnWidth = 17.330729166666668
$choices.Columns.Item(2).ColumnWidth = 19.666666666666668

# Re-select E2 on properties last so it remains the active sheet/tab, matching
# the workbook's activeTab / tabSelected state.
$properties.Range("E2").Select()
